# The "steady" worksheet has a helper column (D) that tags each row with a
# transform label ("Diff-" / "Rate/"). Rows 26-54 hold the model's parameter
# block (alpha, beta, gamma, ... ttrend) which no longer carries a Form
# label, so that column is cleared out for those rows; row 55 keeps its
# label and is left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D26:D54").ClearContents()
